# Scheduled-runner profit data refresh for Anima_Profits.xlsx
# Updates computed leve-profit columns (H-N) on several sheets with
# freshly pulled market-board pricing figures.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 4043.5715
$ws.Range("J121").Value = 4350.8335
$ws.Range("L121").Value = 13052.5005
$ws.Range("N121").Value = -16546.5005

$ws.Range("H133").Value = 94530
$ws.Range("J133").Value = 94530
$ws.Range("L133").Value = 94530
$ws.Range("N133").Value = -104650

$ws.Range("H135").Value = 3320.5
$ws.Range("I135").Value = 3320.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 29884.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -27349.5
$ws.Range("N135").ClearContents()

$ws.Range("H138").Value = 2336.2246
$ws.Range("I138").Value = 3081.3845
$ws.Range("J138").Value = 2067.139
$ws.Range("K138").Value = 9244.1535
$ws.Range("L138").Value = 6201.417
$ws.Range("M138").Value = -4104.1535
$ws.Range("N138").Value = -16481.417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2440.1667
$ws.Range("I45").Value = 1432.1428
$ws.Range("J45").Value = 3851.4
$ws.Range("K45").Value = 1432.1428
$ws.Range("L45").Value = 3851.4
$ws.Range("M45").Value = -1055.1428
$ws.Range("N45").Value = -4605.4

$ws.Range("H74").Value = 1333.5416
$ws.Range("I74").Value = 1210.1052
$ws.Range("J74").Value = 1802.6
$ws.Range("K74").Value = 1210.1052
$ws.Range("L74").Value = 1802.6
$ws.Range("M74").Value = -336.1052
$ws.Range("N74").Value = -3550.6

$ws.Range("H77").Value = 1333.5416
$ws.Range("I77").Value = 1210.1052
$ws.Range("J77").Value = 1802.6
$ws.Range("K77").Value = 6050.526
$ws.Range("L77").Value = 9013
$ws.Range("M77").Value = -1682.526
$ws.Range("N77").Value = -17749

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1249.2142
$ws.Range("I22").Value = 1268.3846
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1268.3846
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1095.3846
$ws.Range("N22").Value = -1346

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4254.6577
$ws.Range("I31").Value = 945.9706
$ws.Range("J31").Value = 32378.5
$ws.Range("K31").Value = 945.9706
$ws.Range("L31").Value = 32378.5
$ws.Range("M31").Value = -650.9706
$ws.Range("N31").Value = -32968.5

$ws.Range("H34").Value = 4254.6577
$ws.Range("I34").Value = 945.9706
$ws.Range("J34").Value = 32378.5
$ws.Range("K34").Value = 945.9706
$ws.Range("L34").Value = 32378.5
$ws.Range("M34").Value = -743.9706
$ws.Range("N34").Value = -32782.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 13980.667
$ws.Range("I3").Value = 2030
$ws.Range("J3").Value = 19956
$ws.Range("K3").Value = 6090
$ws.Range("L3").Value = 59868
$ws.Range("M3").Value = -5978
$ws.Range("N3").Value = -60092

$ws.Range("H5").Value = 1803.7142
$ws.Range("I5").Value = 749.75
$ws.Range("J5").Value = 2051.7058
$ws.Range("K5").Value = 2249.25
$ws.Range("L5").Value = 6155.117400000001
$ws.Range("M5").Value = -2137.25
$ws.Range("N5").Value = -6379.117400000001

$ws.Range("H34").Value = 12195507
$ws.Range("J34").Value = 14706331
$ws.Range("L34").Value = 44118993
$ws.Range("N34").Value = -44119161

$ws.Range("H70").Value = 2150
$ws.Range("I70").Value = 300
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 900
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -585
$ws.Range("N70").Value = -12630

$ws.Range("H73").Value = 2150
$ws.Range("I73").Value = 300
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 900
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = 192
$ws.Range("N73").Value = -14184

$ws.Range("H88").Value = 6861.7144
$ws.Range("J88").Value = 6861.7144
$ws.Range("L88").Value = 20585.1432
$ws.Range("N88").Value = -21441.1432

$ws.Range("H91").Value = 6861.7144
$ws.Range("J91").Value = 6861.7144
$ws.Range("L91").Value = 20585.1432
$ws.Range("N91").Value = -23549.1432

$ws.Range("H94").Value = 2800
$ws.Range("J94").Value = 3200
$ws.Range("L94").Value = 9600
$ws.Range("N94").Value = -10952

$ws.Range("H100").Value = 2663.3333
$ws.Range("J100").Value = 2663.3333
$ws.Range("L100").Value = 7989.999899999999
$ws.Range("N100").Value = -9611.999899999999

$ws.Range("H103").Value = 1945.7142
$ws.Range("I103").Value = 1040
$ws.Range("J103").Value = 2625
$ws.Range("K103").Value = 3120
$ws.Range("L103").Value = 7875
$ws.Range("M103").Value = -2241
$ws.Range("N103").Value = -9633

$ws.Range("H106").Value = 7528.7144
$ws.Range("I106").Value = 1501
$ws.Range("J106").Value = 8533.333
$ws.Range("K106").Value = 4503
$ws.Range("L106").Value = 25599.999
$ws.Range("M106").Value = -3557
$ws.Range("N106").Value = -27491.999

$ws.Range("H117").Value = 2032
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 2032
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 6096
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -12980

$ws.Range("H135").Value = 1803.7142
$ws.Range("I135").Value = 749.75
$ws.Range("J135").Value = 2051.7058
$ws.Range("K135").Value = 6747.75
$ws.Range("L135").Value = 18465.3522
$ws.Range("M135").Value = -4212.75
$ws.Range("N135").Value = -23535.3522

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 35651.332
$ws.Range("J10").Value = 35651.332
$ws.Range("L10").Value = 35651.332
$ws.Range("N10").Value = -35931.332

$ws.Range("H22").Value = 13406.333
$ws.Range("I22").Value = 1700
$ws.Range("J22").Value = 16751
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 16751
$ws.Range("M22").Value = -1405
$ws.Range("N22").Value = -17341

$ws.Range("H27").Value = 13406.333
$ws.Range("I27").Value = 1700
$ws.Range("J27").Value = 16751
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 16751
$ws.Range("M27").Value = -1593
$ws.Range("N27").Value = -16965

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 67504.5
$ws.Range("J10").Value = 67504.5
$ws.Range("L10").Value = 67504.5
$ws.Range("N10").Value = -67842.5

$ws.Range("H132").Value = 6175791
$ws.Range("I132").Value = 3408.2
$ws.Range("J132").Value = 13891270
$ws.Range("K132").Value = 10224.6
$ws.Range("L132").Value = 41673810
$ws.Range("M132").Value = -7694.599999999999
$ws.Range("N132").Value = -41678870
